# Generate Report for Handoff
# Replace the old handoff-file identifier (UUID-based base name) and related
# timestamps with the new ones produced by the latest handoff generation run.

$oldId = "99c60e63-25c8-4590-a1ed-e5d74c9fea5a"
$newId = "25aab421-fdc5-4650-8d89-a4bf30dd8e1c"

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------

# File Name (A2) - plain text value, no hyperlink attached.
$wsOverview.Range("A2").Value = "$newId.md"

# Path And Name (B2) - has a hyperlink whose display text must change too,
# while the link target (pointing at the immutable commit blob) stays the
# same. Recreate the hyperlink on the cell with the updated display text.
$hlTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81299ec43aed6428c168bfedabdaa1538f948f95/e2e/$oldId.md"
$rngB2 = $wsOverview.Range("B2")
$rngB2.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($rngB2, $hlTarget, "", "", "e2e\$newId.md")

# Latest HO Xliff Generate Date (G2)
$wsOverview.Range("G2").Value = "2016-09-03 17:05:43"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------

# Source File Name (A2) - hyperlink, same target, new display text.
$rngZhA2 = $wsZhCn.Range("A2")
$rngZhA2.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($rngZhA2, $hlTarget, "", "", "$newId.md")

# Latest Handoff File (G2)
$wsZhCn.Range("G2").Value = "$newId.8c72e73b03fa6c889b7843ff55873b3b7372c05f.zh-cn.xlf"

# Latest Handoff Datetime (H2)
$wsZhCn.Range("H2").Value = "2016-09-03 17:05:38"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------

# Source File Name (A2) - hyperlink, same target, new display text.
$rngDeA2 = $wsDeDe.Range("A2")
$rngDeA2.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($rngDeA2, $hlTarget, "", "", "$newId.md")

# Latest Handoff File (G2)
$wsDeDe.Range("G2").Value = "$newId.8c72e73b03fa6c889b7843ff55873b3b7372c05f.de-de.xlf"

# Latest Handoff Datetime (H2)
$wsDeDe.Range("H2").Value = "2016-09-03 17:05:43"
